# Update Daily Report: 2026-02-20
# Appends the new day's raw records (date serial 46072) to Daily_Data,
# then refreshes the dependent Today_Summary and Monthly_Stats rollups
# to reflect that day's figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Daily_Data: append 24 new rows (one per depository/region type)
#    for the new reporting date.
# ---------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

$newDate = 46072

# Columns: B=Region_Type, C=PREV_TOTAL, D=RECEIVED, E=WITHDRAWN, F=NET_CHANGE, G=ADJUSTMENT, H=TOTAL_TODAY
$newRows = @(
    @("ASAHI DEPOSITORY LLC Registered", 23301775.992, 0, 0, 0, 0, 23301775.992),
    @("ASAHI DEPOSITORY LLC Eligible", 2748893.808, 0, 0, 0, 0, 2748893.808),
    @("BRINK'S, INC. Registered", 15782712.636, 0, 0, 0, 0, 15782712.636),
    @("BRINK'S, INC. Eligible", 39336942.517, 0, 0, 0, 0, 39336942.517),
    @("CNT DEPOSITORY, INC. Registered", 12174851.569, 0, 0, 0, 0, 12174851.569),
    @("CNT DEPOSITORY, INC. Eligible", 14018899.428, 0, 0, 0, 0, 14018899.428),
    @("DELAWARE DEPOSITORY Registered", 1532776.423, 0, 0, 0, 0, 1532776.423),
    @("DELAWARE DEPOSITORY Eligible", 16276669.15, 0, 4000.8, -4000.8, -39.5, 16272628.85),
    @("HSBC BANK, USA Registered", 3412157.57, 0, 0, 0, 0, 3412157.57),
    @("HSBC BANK, USA Eligible", 20264362.853, 0, 625953.5, -625953.5, 0, 19638409.353),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 273789.87, 0, 0, 0, 0, 273789.87),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 3642206.244, 0, 346959.6, -346959.6, 0, 3295246.644),
    @("JP MORGAN CHASE BANK NA Registered", 12000343.77, 0, 0, 0, 0, 12000343.77),
    @("JP MORGAN CHASE BANK NA Eligible", 143854408.433, 0, 0, 0, 0, 143854408.433),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 6311885.937, 0, 0, 0, 0, 6311885.937),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 24033585.186, 0, 0, 0, 0, 24033585.186),
    @("MALCA-AMIT ARMORED, INC. Registered", 0, 0, 0, 0, 0, 0),
    @("MALCA-AMIT ARMORED, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("MALCA-AMIT USA, LLC Registered", 949634.064, 0, 0, 0, 0, 949634.064),
    @("MALCA-AMIT USA, LLC Eligible", 1073898.377, 0, 0, 0, 0, 1073898.377),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 6219630.033, 0, 0, 0, 0, 6219630.033),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 12448651.307, 0, 192635.4, -192635.4, 0, 12256015.907),
    @("STONEX PRECIOUS METALS LLC Registered", 6231501.4, 0, 0, 0, 0, 6231501.4),
    @("STONEX PRECIOUS METALS LLC Eligible", 1537051.72, 0, 0, 0, 0, 1537051.72)
)

$startRow = 194
$r = $startRow
foreach ($row in $newRows) {
    $daily.Cells.Item($r, 1).Value = $newDate
    $daily.Cells.Item($r, 2).Value = $row[0]
    $daily.Cells.Item($r, 3).Value = $row[1]
    $daily.Cells.Item($r, 4).Value = $row[2]
    $daily.Cells.Item($r, 5).Value = $row[3]
    $daily.Cells.Item($r, 6).Value = $row[4]
    $daily.Cells.Item($r, 7).Value = $row[5]
    $daily.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
$endRow = $r - 1

# Match the Date column's existing date/time display format so the new
# cells pick up the same style as the rest of column A.
$daily.Range("A" + $startRow + ":A" + $endRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------
# 2) Today_Summary: per-depository Eligible/Registered/Total_Stock,
#    refreshed with the new day's Eligible totals (Registered totals
#    for these four depositories are unchanged).
# ---------------------------------------------------------------
$today = $wb.Worksheets.Item("Today_Summary")

$today.Cells.Item(5, 2).Value = 16272628.85
$today.Cells.Item(5, 4).Value = 17805405.273

$today.Cells.Item(6, 2).Value = 19638409.353
$today.Cells.Item(6, 4).Value = 23050566.923

$today.Cells.Item(7, 2).Value = 3295246.644
$today.Cells.Item(7, 4).Value = 3569036.514

$today.Cells.Item(12, 2).Value = 12256015.907
$today.Cells.Item(12, 4).Value = 18475645.94

# ---------------------------------------------------------------
# 3) Monthly_Stats: grand total row + per-depository monthly
#    RECEIVED/WITHDRAWN/TOTAL_TODAY detail rows.
# ---------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")

$monthly.Cells.Item(2, 2).Value = 278065980.223
$monthly.Cells.Item(2, 4).Value = 366257039.487

# DELAWARE DEPOSITORY Eligible (row 13): WITHDRAWN += 4000.8, TOTAL_TODAY -> new total
$monthly.Cells.Item(13, 4).Value = 190667.219
$monthly.Cells.Item(13, 5).Value = 16272628.85

# HSBC BANK, USA Eligible (row 15): WITHDRAWN += 625953.5, TOTAL_TODAY -> new total
$monthly.Cells.Item(15, 4).Value = 1682646.81
$monthly.Cells.Item(15, 5).Value = 19638409.353

# INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible (row 17): WITHDRAWN += 346959.6, TOTAL_TODAY -> new total
$monthly.Cells.Item(17, 4).Value = 346959.6
$monthly.Cells.Item(17, 5).Value = 3295246.644

# MANFRA, TORDELLA & BROOKES, LLC Eligible (row 27): WITHDRAWN += 192635.4, TOTAL_TODAY -> new total
$monthly.Cells.Item(27, 4).Value = 1042561.223
$monthly.Cells.Item(27, 5).Value = 12256015.907
